# Natmi following Dr Hou advice
# Adds "ECs" (endothelial cells) as a sending cluster, completing the full
# ECs / FAPs / sCs sender-target cross product for the Sema3c -> Nrp1 edge,
# and refreshes all derived NATMI statistics (rows 2-10, columns E-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema3c", "Nrp1", "ECs", 1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 24.30972217721378, 218.787499594924, 0.003227476896623221, 0.003227476896623221),
    @("ECs", "Sema3c", "Nrp1", "FAPs", 1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 15.05043952253111, 135.45395570278, 0.001998169517894569, 0.001998169517894569),
    @("ECs", "Sema3c", "Nrp1", "sCs", 1, 0.3333333333333333, 0.2364093333333333, 0.709228, 0.005805733041453686, 0.005805733041453687, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 4.369278291126666, 39.32350462014, 0.0005800866269358966, 0.0005800866269358966),
    @("FAPs", "Sema3c", "Nrp1", "ECs", 3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 4020.272525022293, 36182.45272520064, 0.533750924755558, 0.533750924755558),
    @("FAPs", "Sema3c", "Nrp1", "FAPs", 3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 2488.998766043338, 22400.98889439004, 0.3304515763104163, 0.3304515763104163),
    @("FAPs", "Sema3c", "Nrp1", "sCs", 3, 1, 39.09670133333334, 117.290104, 0.9601355730855794, 0.9601355730855795, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 722.5787830869467, 6503.20904778252, 0.09593307201960516, 0.09593307201960516),
    @("sCs", "Sema3c", "Nrp1", "ECs", 3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 102.8289443333334, 308.486833, 0.5559120396302444, 0.5559120396302443, 142.6103094749405, 1283.492785274464, 0.01893363797806311, 0.01893363797806311),
    @("sCs", "Sema3c", "Nrp1", "FAPs", 3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 63.66262833333334, 190.987885, 0.3441717873742006, 0.3441717873742006, 88.29174691489779, 794.62572223408, 0.01172204154588974, 0.01172204154588974),
    @("sCs", "Sema3c", "Nrp1", "sCs", 3, 1, 1.386869333333333, 4.160608, 0.03405869387296686, 0.03405869387296687, 3, 1, 18.481835, 55.445505, 0.09991617299555507, 0.09991617299555505, 25.63189018522667, 230.68701166704, 0.003403014349014008, 0.003403014349014008)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

